$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Content.Find.Execute("2025-07-22 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-07-23 Wednesday", 2)

# Update the division problems in the table. Addressing by row/column
# avoids ambiguity where the same problem text (e.g. "71÷2=") appears
# more than once in the document.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "54÷8="
$t.Cell(1, 2).Range.Text = "20÷4="
$t.Cell(1, 3).Range.Text = "32÷6="
$t.Cell(1, 4).Range.Text = "63÷3="
$t.Cell(1, 5).Range.Text = "89÷7="

$t.Cell(5, 1).Range.Text = "87÷4="
$t.Cell(5, 2).Range.Text = "75÷3="
$t.Cell(5, 3).Range.Text = "33÷3="
$t.Cell(5, 4).Range.Text = "68÷7="
$t.Cell(5, 5).Range.Text = "13÷6="

$t.Cell(9, 1).Range.Text = "97÷3="
$t.Cell(9, 2).Range.Text = "70÷5="
$t.Cell(9, 3).Range.Text = "64÷4="
$t.Cell(9, 4).Range.Text = "77÷4="
$t.Cell(9, 5).Range.Text = "91÷5="

$t.Cell(13, 1).Range.Text = "56÷4="
$t.Cell(13, 2).Range.Text = "70÷2="
$t.Cell(13, 3).Range.Text = "31÷2="
$t.Cell(13, 4).Range.Text = "41÷4="
$t.Cell(13, 5).Range.Text = "93÷2="

$t.Cell(17, 1).Range.Text = "87÷2="
$t.Cell(17, 2).Range.Text = "89÷8="
$t.Cell(17, 3).Range.Text = "60÷9="
$t.Cell(17, 4).Range.Text = "75÷9="
$t.Cell(17, 5).Range.Text = "42÷9="
